$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8053624310004466
$ws.Range("C2").Value = 0.2430607420111528
$ws.Range("D2").Value = 0.07976951290095258
$ws.Range("E2").Value = 0.4261583003367519
$ws.Range("G2").Value = 0.2829097545395314
$ws.Range("H2").Value = 0.4154558441493634
$ws.Range("I2").Value = 0.2692619540102292
$ws.Range("N2").Value = 0.7912964530252466
$ws.Range("O2").Value = 1.319573005448774

$ws.Range("B3").Value = 0.7040990595721155
$ws.Range("C3").Value = 0.2141160009225587
$ws.Range("D3").Value = 0.07218912734981586
$ws.Range("E3").Value = 0.3716950293549672
$ws.Range("G3").Value = 0.2735515995547644
$ws.Range("H3").Value = 0.4153253780719126
$ws.Range("I3").Value = 0.2721510507137133
$ws.Range("N3").Value = 0.7924410893878928
$ws.Range("O3").Value = 1.299146145657829

$ws.Range("B4").Value = 0.641768240904014
$ws.Range("C4").Value = 0.1962558766066422
$ws.Range("D4").Value = 0.067569147290925
$ws.Range("E4").Value = 0.3383453989565481
$ws.Range("G4").Value = 0.268097865472825
$ws.Range("H4").Value = 0.4155251206244088
$ws.Range("I4").Value = 0.2741748367591264
$ws.Range("N4").Value = 0.7934924659542872
$ws.Range("O4").Value = 1.287798132332171

$ws.Range("B5").Value = 0.6163303107064166
$ws.Range("C5").Value = 0.1889558737278207
$ws.Range("D5").Value = 0.06569510267496526
$ws.Range("E5").Value = 0.3247762222429884
$ws.Range("G5").Value = 0.2659484127097045
$ws.Range("H5").Value = 0.4156767555073628
$ws.Range("I5").Value = 0.2750621419866732
$ws.Range("N5").Value = 0.7940087357526551
$ws.Range("O5").Value = 1.283472681428663

$ws.Range("B6").Value = 0.6121041278771884
$ws.Range("C6").Value = 0.1877424012540985
$ws.Range("D6").Value = 0.0653844400210204
$ws.Range("E6").Value = 0.3225242901737317
$ws.Range("G6").Value = 0.2655958916806611
$ws.Range("H6").Value = 0.4157061727293865
$ws.Range("I6").Value = 0.2752132520942361
$ws.Range("N6").Value = 0.7940997716770326
$ws.Range("O6").Value = 1.282772461159368

$ws.Range("B7").Value = 0.6414253266770515
$ws.Range("C7").Value = 0.1961575143731125
$ws.Range("D7").Value = 0.067543838329712
$ws.Range("E7").Value = 0.3381623173800676
$ws.Range("G7").Value = 0.2680685823256397
$ws.Range("H7").Value = 0.4155268814108837
$ws.Range("I7").Value = 0.2741865501327077
$ws.Range("N7").Value = 0.793499072713594
$ws.Range("O7").Value = 1.287738589086189

$ws.Range("B8").Value = 0.7704796668114682
$ws.Range("C8").Value = 0.2330989905123886
$ws.Range("D8").Value = 0.07714864193556537
$ws.Range("E8").Value = 0.4073593540735629
$ws.Range("G8").Value = 0.279622115344651
$ws.Range("H8").Value = 0.4153526956168321
$ws.Range("I8").Value = 0.2702061123012349
$ws.Range("N8").Value = 0.7916188808979925
$ws.Range("O8").Value = 1.312281167607694

$ws.Range("B9").Value = 1.022283190993846
$ws.Range("C9").Value = 0.3048356044010347
$ws.Range("D9").Value = 0.09625843345085627
$ws.Range("E9").Value = 0.5438654983744726
$ws.Range("G9").Value = 0.3046208002748614
$ws.Range("H9").Value = 0.4172382159653552
$ws.Range("I9").Value = 0.2643938076348746
$ws.Range("N9").Value = 0.7906908200730953
$ws.Range("O9").Value = 1.369948641587769

$ws.Range("B10").Value = 1.206466752585186
$ws.Range("C10").Value = 0.3571050245740537
$ws.Range("D10").Value = 0.1104698725585678
$ws.Range("E10").Value = 0.6447828074138613
$ws.Range("G10").Value = 0.3244502324685499
$ws.Range("H10").Value = 0.4199913196834331
$ws.Range("I10").Value = 0.2613536708444606
$ws.Range("N10").Value = 0.7916833697626657
$ws.Range("O10").Value = 1.418228385698455

$ws.Range("B11").Value = 1.290071784059194
$ws.Range("C11").Value = 0.380788379983187
$ws.Range("D11").Value = 0.11697319245269
$ws.Range("E11").Value = 0.6908578334543591
$ws.Range("G11").Value = 0.3337964840977179
$ws.Range("H11").Value = 0.4215430416844583
$ws.Range("I11").Value = 0.2602409959483651
$ws.Range("N11").Value = 0.7924970884570683
$ws.Range("O11").Value = 1.441496808335188

$ws.Range("B12").Value = 1.321703772474734
$ws.Range("C12").Value = 0.389742914226872
$ws.Range("D12").Value = 0.1194414037149585
$ws.Range("E12").Value = 0.7083314840976556
$ws.Range("G12").Value = 0.3373830824541244
$ws.Range("H12").Value = 0.4221738487475761
$ws.Range("I12").Value = 0.2598587874370644
$ws.Range("N12").Value = 0.792857182114787
$ws.Range("O12").Value = 1.450497294622693

$ws.Range("B13").Value = 1.314892499562006
$ws.Range("C13").Value = 0.3878150150122792
$ws.Range("D13").Value = 0.1189095834145348
$ws.Range("E13").Value = 0.7045670308366567
$ws.Range("G13").Value = 0.3366085287594842
$ws.Range("H13").Value = 0.4220360690814573
$ws.Range("I13").Value = 0.2599393579493814
$ws.Range("N13").Value = 0.7927773207179172
$ws.Range("O13").Value = 1.448550438788999

$ws.Range("B14").Value = 1.292674724338383
$ws.Range("C14").Value = 0.3815253538617753
$ws.Range("D14").Value = 0.117176142465766
$ws.Range("E14").Value = 0.6922948678940912
$ws.Range("G14").Value = 0.3340906031593391
$ws.Range("H14").Value = 0.4215940717668758
$ws.Range("I14").Value = 0.2602087655124343
$ws.Range("N14").Value = 0.7925256730147794
$ws.Range("O14").Value = 1.442233482632076

$ws.Range("B15").Value = 1.279062082978839
$ws.Range("C15").Value = 0.3776709442585116
$ws.Range("D15").Value = 0.1161150823874522
$ws.Range("E15").Value = 0.6847812640872917
$ws.Range("G15").Value = 0.3325544866149386
$ws.Range("H15").Value = 0.4213289668501119
$ws.Range("I15").Value = 0.2603788902444677
$ws.Range("N15").Value = 0.7923782940422655
$ws.Range("O15").Value = 1.438388853119591

$ws.Range("B16").Value = 1.200999205588005
$ws.Range("C16").Value = 0.3555553404717102
$ws.Range("D16").Value = 0.1100456405173134
$ws.Range("E16").Value = 0.6417752309760374
$ws.Range("G16").Value = 0.3238460333736697
$ws.Range("H16").Value = 0.4198959478975297
$ws.Range("I16").Value = 0.2614318488204255
$ws.Range("N16").Value = 0.7916374686683696
$ws.Range("O16").Value = 1.416734140311974

$ws.Range("B17").Value = 1.153062806073478
$ws.Range("C17").Value = 0.3419637878972424
$ws.Range("D17").Value = 0.1063320962832819
$ws.Range("E17").Value = 0.6154366654292716
$ws.Range("G17").Value = 0.3185874898250063
$ws.Range("H17").Value = 0.4190936151605342
$ws.Range("I17").Value = 0.2621472257548625
$ws.Range("N17").Value = 0.791275671847302
$ws.Range("O17").Value = 1.403785152989343

$ws.Range("B18").Value = 1.125474103808187
$ws.Range("C18").Value = 0.3341374303607552
$ws.Range("D18").Value = 0.1041997835428106
$ws.Range("E18").Value = 0.6003030401099494
$ws.Range("G18").Value = 0.3155935431691006
$ws.Range("H18").Value = 0.4186602983889998
$ws.Range("I18").Value = 0.2625841173189016
$ws.Range("N18").Value = 0.7911016753837998
$ws.Range("O18").Value = 1.396460050301471

$ws.Range("B19").Value = 1.116130174204727
$ws.Range("C19").Value = 0.3314860478143089
$ws.Range("D19").Value = 0.1034784403433662
$ws.Range("E19").Value = 0.5951816771887337
$ws.Range("G19").Value = 0.3145850884941979
$ws.Range("H19").Value = 0.4185184170342211
$ws.Range("I19").Value = 0.2627364000042292
$ws.Range("N19").Value = 0.7910486238095302
$ws.Range("O19").Value = 1.394000944585144

$ws.Range("B20").Value = 1.158167482357612
$ws.Range("C20").Value = 0.3434115504925614
$ws.Range("D20").Value = 0.1067270346829616
$ws.Range("E20").Value = 0.6182388196532997
$ws.Range("G20").Value = 0.3191440969600592
$ws.Range("H20").Value = 0.4191761088335397
$ws.Range("I20").Value = 0.2620684392508892
$ws.Range("N20").Value = 0.7913106575127244
$ws.Range("O20").Value = 1.405150873481546

$ws.Range("B21").Value = 1.299201378811574
$ws.Range("C21").Value = 0.3833731576210653
$ws.Range("D21").Value = 0.1176851454873997
$ws.Range("E21").Value = 0.6958987770053966
$ws.Range("G21").Value = 0.3348288889173432
$ws.Range("H21").Value = 0.4217227232650345
$ws.Range("I21").Value = 0.2601285696601039
$ws.Range("N21").Value = 0.7925981789636865
$ws.Range("O21").Value = 1.444083778068062

$ws.Range("B22").Value = 1.391214701534636
$ws.Range("C22").Value = 0.4094095871954551
$ws.Range("D22").Value = 0.1248792321701728
$ws.Range("E22").Value = 0.7468067570364241
$ws.Range("G22").Value = 0.3453561977830759
$ws.Range("H22").Value = 0.4236389825825313
$ws.Range("I22").Value = 0.2590890001771697
$ws.Range("N22").Value = 0.7937424332167495
$ws.Range("O22").Value = 1.470632404613838

$ws.Range("B23").Value = 1.342120602643718
$ws.Range("C23").Value = 0.3955209400106696
$ws.Range("D23").Value = 0.1210366501441342
$ws.Range("E23").Value = 0.7196215744520345
$ws.Range("G23").Value = 0.3397121083495449
$ws.Range("H23").Value = 0.4225931364066184
$ws.Range("I23").Value = 0.2596228643692662
$ws.Range("N23").Value = 0.7931040591310818
$ws.Range("O23").Value = 1.456361422674149

$ws.Range("B24").Value = 1.155859749770116
$ws.Range("C24").Value = 0.342757055609809
$ws.Range("D24").Value = 0.1065484747707899
$ws.Range("E24").Value = 0.6169719386100354
$ws.Range("G24").Value = 0.3188923637701038
$ws.Range("H24").Value = 0.419138726372708
$ws.Range("I24").Value = 0.2621039788674047
$ws.Range("N24").Value = 0.7912947345538015
$ws.Range("O24").Value = 1.404533059290458

$ws.Range("B25").Value = 0.9543040993675618
$ws.Range("C25").Value = 0.2855050429776043
$ws.Range("D25").Value = 0.09105890053938026
$ws.Range("E25").Value = 0.5068365844614391
$ws.Range("G25").Value = 0.2976036262695487
$ws.Range("H25").Value = 0.4164886654822766
$ws.Range("I25").Value = 0.2657512451560784
$ws.Range("N25").Value = 0.7906473227258317
$ws.Range("O25").Value = 1.3533168933985
